$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row data to append: each entry describes the text for column A and B
# of a new row, plus which "style" (1 = plain existing body style,
# 2 = new Google-Sans / white-fill style) each cell should use.
# ---------------------------------------------------------------------
$rowData = @(
    @{ Row = 116; SA = "1"; TA = "China's CO2 emissions from energy and cement"; SB = "1"; TB = "中国化石燃料和水泥的二氧化碳排放量" },
    @{ Row = 117; SA = "1"; TA = "Mt CO2 / quarter"; SB = "1"; TB = "亿吨二氧化碳/季度" },
    @{ Row = 118; SA = "1"; TA = "Trends in China CO2 by fuel and sector"; SB = "1"; TB = "不同行业化石燃料和水泥的二氧化碳排放量变化趋势" },
    @{ Row = 119; SA = "1"; TA = "pre-COVID trendline"; SB = "1"; TB = "疫情前变化趋势" },
    @{ Row = 120; SA = "1"; TA = "First COVID-19 lockdown"; SB = "1"; TB = "新冠疫情爆发" },
    @{ Row = 121; SA = "1"; TA = "Mt/year, 12 month moving sum"; SB = "1"; TB = "亿吨/年，12个月移动平均值" },
    @{ Row = 122; SA = "1"; TA = "Cement"; SB = "1"; TB = "水泥" },
    @{ Row = 123; SA = "1"; TA = "Coking Coal"; SB = "1"; TB = "炼焦煤" },
    @{ Row = 124; SA = "1"; TA = "Natural Gas"; SB = "1"; TB = "天然气" },
    @{ Row = 125; SA = "1"; TA = "Oil Products"; SB = "1"; TB = "石油产品" },
    @{ Row = 126; SA = "1"; TA = "Steam coal"; SB = "1"; TB = "动力煤" },
    @{ Row = 127; SA = "1"; TA = "Steam coal"; SB = "1"; TB = "动力煤" },
    @{ Row = 128; SA = "1"; TA = "Automobiles"; SB = "1"; TB = "汽车" },
    @{ Row = 129; SA = "1"; TA = "New Energy Vehicles"; SB = "1"; TB = "新能源汽车" },
    @{ Row = 130; SA = "1"; TA = "COVID-19 lockdown"; SB = "1"; TB = "新冠疫情爆发" },
    @{ Row = 131; SA = "1"; TA = "Vehicles Production"; SB = "1"; TB = "乘用车产量" },
    @{ Row = 132; SA = "1"; TA = "Million units, 12-month moving sum"; SB = "1"; TB = "百万辆，12个月移动平均值" },
    @{ Row = 133; SA = "1"; TA = "Automobiles"; SB = "1"; TB = "汽车" },
    @{ Row = 134; SA = "1"; TA = "New Energy Vehicles"; SB = "1"; TB = "新能源汽车" },
    @{ Row = 135; SA = "1"; TA = "cumulative sales over 10 years"; SB = "1"; TB = "过去十年总销量" },
    @{ Row = 136; SA = "1"; TA = "new sales, 3-month mean"; SB = "1"; TB = "新车总销量，三个月平均值" },
    @{ Row = 137; SA = "1"; TA = "new energy vehicle share"; SB = "1"; TB = "新能源车占比" },
    @{ Row = 138; SA = "1"; TA = "Power Sector Coal Consumption in China"; SB = "1"; TB = "电力行业煤炭消耗量" },
    @{ Row = 139; SA = "1"; TA = "Mt/year, 12-month moving sum"; SB = "1"; TB = "亿吨/年，12个月移动平均值" },
    @{ Row = 140; SA = "1"; TA = "predicted based on output"; SB = "1"; TB = "根据产量测算" },
    @{ Row = 141; SA = "1"; TA = "reported"; SB = "1"; TB = "公布的数据" },
    @{ Row = 142; SA = "1"; TA = "without drought&heatwave"; SB = "1"; TB = "剔除干旱和水电产量下降因素" },
    @{ Row = 143; SA = "2"; TA = "All Sectors"; SB = "2"; TB = "所有行业" },
    @{ Row = 144; SA = "1"; TA = "Non-power use"; SB = "2"; TB = "非电力用途" },
    @{ Row = 145; SA = "1"; TA = "Power Industry"; SB = "2"; TB = "电力行业" },
    @{ Row = 146; SA = "2"; TA = "Metallurgy Industry"; SB = "2"; TB = "冶金行业" },
    @{ Row = 147; SA = "1"; TA = "Quarterly"; SB = "1"; TB = "每季度" },
    @{ Row = 148; SA = "1"; TA = "Mt/year"; SB = "2"; TB = "亿吨/年" },
    @{ Row = 149; SA = "1"; TA = "Power"; SB = "2"; TB = "电力行业" },
    @{ Row = 150; SA = "1"; TA = "Total"; SB = "2"; TB = "所有行业" },
    @{ Row = 151; SA = "1"; TA = "Coal consumption in China"; SB = "1"; TB = "煤炭消耗量" }
)

foreach ($item in $rowData) {
    $ws.Cells.Item($item.Row, 1).Value = $item.TA
    $ws.Cells.Item($item.Row, 2).Value = $item.TB
}

# ---------------------------------------------------------------------
# Apply the "plain" body style (same formatting as the existing rows,
# e.g. row 115) to every new cell first, by copying the format from an
# existing style-1 cell. This reuses the existing cellXf instead of
# creating a new one.
# ---------------------------------------------------------------------
$ws.Range("A115:B115").Copy() | Out-Null
foreach ($item in $rowData) {
    $ws.Range("A$($item.Row):B$($item.Row)").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Build the new "Google Sans" style (font size 9, color #1F1F1F, solid
# white fill) on a scratch cell, then propagate it with copy/paste so
# that all target cells share a single new cellXf entry.
# ---------------------------------------------------------------------
$template = $ws.Cells.Item(300, 10)
$template.Value = "style-template"
$template.Font.Name = '"Google Sans"'
$template.Font.Size = 9
$template.Font.Color = 2039583
$template.Interior.Color = 16777215
$template.Interior.PatternColor = 16777215

$template.Copy() | Out-Null
foreach ($item in $rowData) {
    if ($item.SA -eq "2") {
        $ws.Cells.Item($item.Row, 1).PasteSpecial(-4122) | Out-Null
    }
    if ($item.SB -eq "2") {
        $ws.Cells.Item($item.Row, 2).PasteSpecial(-4122) | Out-Null
    }
}
$excel.CutCopyMode = 0

$template.Clear() | Out-Null
